$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: spreadsheet row number, then new text for columns D (Price),
# E (Volume 1h) and G (Hora). A column is omitted when the diff left it unchanged.
$updates = @(
    @{Row=2; D='246.16'; E='-0.31%'; G='10'},
    @{Row=3; D='30.19'; E='0.03%'; G='10'},
    @{Row=4; D='5.153'; E='-0.49%'; G='10'},
    @{Row=5; D='0.05773'; E='0.59%'; G='10'},
    @{Row=6; D='6.655'; E='0.93%'; G='10'},
    @{Row=7; D='3.244'; E='6.27%'; G='10'},
    @{Row=8; D='0.8492'; E='-1.12%'; G='10'},
    @{Row=9; D='0.8555'; E='-2.21%'; G='10'},
    @{Row=10; D='0.1390'; E='2.23%'; G='10'},
    @{Row=11; D='0.07084'; E='-0.01%'; G='10'},
    @{Row=12; D='0.03259'; E='11.82%'; G='10'},
    @{Row=13; D='0.09374'; E='-0.20%'; G='10'},
    @{Row=14; D='0.001524'; E='1.04%'; G='10'},
    @{Row=15; D='0.0005946'; E='-94.20%'; G='10'},
    @{Row=16; D='0.006054'; E='-0.34%'; G='10'},
    @{Row=17; D='3.523'; E='0.55%'; G='10'},
    @{Row=18; D='2.187'; E='-4.19%'; G='10'},
    @{Row=19; D='0.3165'; E='-0.58%'; G='10'},
    @{Row=20; D='0.03387'; E='3.28%'; G='10'},
    @{Row=21; E='0.25%'; G='10'},
    @{Row=22; D='3.492'; E='-3.16%'; G='10'},
    @{Row=23; E='2.23%'; G='10'},
    @{Row=24; D='0.04108'; E='-0.92%'; G='10'},
    @{Row=25; D='0.001228'; E='1.08%'; G='10'},
    @{Row=26; E='-7.96%'; G='10'},
    @{Row=27; D='0.0001200'; E='1.78%'; G='10'},
    @{Row=28; E='4.24%'; G='10'},
    @{Row=29; G='10'},
    @{Row=30; G='10'},
    @{Row=31; G='10'},
    @{Row=32; G='10'},
    @{Row=33; G='10'},
    @{Row=34; G='10'},
    @{Row=35; G='10'},
    @{Row=36; G='10'},
    @{Row=37; G='10'},
    @{Row=38; G='10'},
    @{Row=39; G='10'},
    @{Row=40; D='0.03760'; E='-0.69%'; G='10'},
    @{Row=41; D='0.1072'; E='0.00%'; G='10'},
    @{Row=42; D='0.002199'; E='0.08%'; G='10'},
    @{Row=43; D='0.002949'; E='-47.74%'; G='10'},
    @{Row=44; D='0.009954'; E='-0.48%'; G='10'},
    @{Row=45; D='0.00005477'; E='7.19%'; G='10'},
    @{Row=46; D='0.00000000750'; E='0.08%'; G='10'},
    @{Row=47; D='0.07098'; E='-20.15%'; G='10'},
    @{Row=48; E='-10.85%'; G='10'},
    @{Row=49; D='0.00002099'; E='0.08%'; G='10'},
    @{Row=50; D='0.0001999'; E='0.08%'; G='10'},
    @{Row=51; G='10'}
)

foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @("D", "E", "G")) {
        if ($u.ContainsKey($col)) {
            $cell = $ws.Range("$col$r")
            $cell.NumberFormat = "@"
            $cell.Value = $u[$col]
        }
    }
}
